# Regenerate save_data to use K (strike count) instead of Strike#,
# recalculated from the regenerated std/mean + s_vals computation.
# Only column G ("K") changes; every other column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,0,1,1,1,2,2,1,0,1,0,0,1,1,0,0,0,1,1,2,0,0,1,1,0,2,1,2,0,0,1,1,3,2,3,2,1,2,0,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
